# Updated cryptos list on Mon Jan  1 17:53:21 UTC 2024 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns with the latest scraped figures.
# Rows 50/51 additionally swap place (THORChain now ranks above FraxShare).
#
# Some new Price values are plain decimals (e.g. "1.00", "311.52") that Excel
# would otherwise auto-parse as numbers; those cells are pre-formatted as
# Text ("@") so the value is stored verbatim as a string, matching the
# source data (which also stores two-dot "thousands" prices like
# "43.116.63" as text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.116.63'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '2.314.12'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.52'
$ws.Range('E5').Value = '  -1.63%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.22'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.10'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0915'
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.39'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.990'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.25'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').Value = '2.664.30'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D17').Value = '2.313.54'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '42.960.63'
$ws.Range('E18').Value = '  +0.63%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.47'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.09'
$ws.Range('E21').Value = '  -13.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.67'
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.50'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '266.90'
$ws.Range('E24').Value = '  +0.22%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.67'
$ws.Range('E27').Value = '  +11.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.06'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.67'
$ws.Range('E30').Value = '  +3.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.47'
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '166.56'
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0877'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('E34').Value = '  +4.94%  '
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.72'
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('E37').Value = '  -1.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0358'
$ws.Range('E38').Value = '  +1.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.81'
$ws.Range('E39').Value = '  +4.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.69'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.61'
$ws.Range('E41').Value = '  +1.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '104.65'
$ws.Range('E42').Value = '  +9.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.91'
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.234'
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.02'
$ws.Range('E45').Value = '  +5.31%  '
$ws.Range('E46').Value = '  +0.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '112.80'
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('D48').Value = '1.657.95'
$ws.Range('E48').Value = '  -3.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '77.16'
$ws.Range('E49').Value = '  -3.63%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.25'
$ws.Range('E50').Value = '  +3.00%  '
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.80'
$ws.Range('E51').Value = '  -0.14%  '
